$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates whose new text is NOT a valid Excel numeric literal
#     (e.g. contains two "." thousand-separators, like "26.405.73") -- plain assignment
#     keeps them stored as text, matching the original inline-string cells. ---
$ws.Range("D2").Value = "26.405.73"
$ws.Range("D3").Value = "1.722.11"
$ws.Range("D10").Value = "1.717.22"
$ws.Range("D17").Value = "26.403.62"
$ws.Range("D21").Value = "1.945.73"

# --- Column D (Price) updates whose new text WOULD be auto-converted into a number by
#     Excel (e.g. "15.48", "0.9997") -- force the cell to Text format first so the value
#     is written verbatim, then restore the default "Normal" style so no stray number
#     formatting is left behind on the cell. ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D5").Value = "242.59"
$ws.Range("D7").Value = "0.4914"
$ws.Range("D8").Value = "0.2612"
$ws.Range("D9").Value = "0.06192"
$ws.Range("D11").Value = "0.07015"
$ws.Range("D12").Value = "15.48"
$ws.Range("D13").Value = "4.560"
$ws.Range("D14").Value = "0.5984"
$ws.Range("D15").Value = "77.18"
$ws.Range("D19").Value = "0.000007146"
$ws.Range("D20").Value = "11.35"
$ws.Range("D22").Value = "4.477"
$ws.Range("D23").Value = "8.575"
$ws.Range("D24").Value = "5.152"
$ws.Range("D25").Value = "137.26"
$ws.Range("D26").Value = "15.21"
$ws.Range("D27").Value = "1.397"
$ws.Range("D28").Value = "106.98"
$ws.Range("D29").Value = "1.701"
$ws.Range("D30").Value = "3.937"
$ws.Range("D31").Value = "0.07951"
$ws.Range("D33").Value = "0.04544"
$ws.Range("D34").Value = "2.602"
$ws.Range("D35").Value = "0.9922"
$ws.Range("D36").Value = "0.6233"
$ws.Range("D37").Value = "0.9234"
$ws.Range("D38").Value = "2.392"
$ws.Range("D40").Value = "0.9997"
$ws.Range("D42").Value = "99.96"
$ws.Range("D43").Value = "5.342"
$ws.Range("D44").Value = "0.3831"
$ws.Range("D46").Value = "0.1162"
$ws.Range("D48").Value = "30.06"
$ws.Range("D49").Value = "7.672"
$ws.Range("D50").Value = "1.233"
$ws.Range("D51").Value = "50.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

# --- Column E (Volume) updates -- the padded "  +x.xx%  " / "  -x.xx%  " text is not a valid
#     Excel numeric/percent literal (leading/trailing spaces), so it is safely assigned as-is.
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  -4.17%  "
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  -6.41%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("E42").Value = "  -3.22%  "
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("E51").Value = "  -0.80%  "
